# "Generate Report for Archive"
#
# The localization status report moved from "Ready for handoff" to
# "In Translation" for the two handed-off files, on every sheet that
# surfaces that status (the Overview rollup columns for each locale, and
# each per-locale detail sheet's Status column). The now-narrower Status
# column no longer needs to be as wide, so it is shrunk to fit.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: per-locale rollup columns (E = zh-cn, F = de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Columns E and F shrink along with the new, shorter status text.
$overview.Range("E1").ColumnWidth = 12.576851254417766
$overview.Range("F1").ColumnWidth = 12.576851254417766

# --- Per-locale detail sheets: Status column (column C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("C1").ColumnWidth = 12.576851254417766

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("C1").ColumnWidth = 12.576851254417766
